# Add the 34 manually-reviewed hospital-unit rows that Excel appended
# to lsh_unit_categories (sheet "lsh_unit_categories"), then restore the
# on-screen selection / active-sheet state recorded by Excel after that edit.

$wb = $excel.ActiveWorkbook
$wsIsolation = $wb.Worksheets.Item("lsh_isolation_categories")
$wsUnits     = $wb.Worksheets.Item("lsh_unit_categories")

# ---- Append rows 148-181 to lsh_unit_categories ----
$wsUnits.Cells.Item(148, 1).Value = "Hb-11F GD Sálfræðiþjónusta"
$wsUnits.Cells.Item(148, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(148, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(148, 4).Value = "home"
$wsUnits.Cells.Item(148, 5).Value = 1
# row flagged for manual hospital-probability review -> reuse existing format from A136
$wsUnits.Range("A136").Copy()
$wsUnits.Range("A148").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsUnits.Cells.Item(149, 1).Value = "Kv-H8 GD Líknardeildar"
$wsUnits.Cells.Item(149, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(149, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(149, 4).Value = "home"
$wsUnits.Cells.Item(149, 5).Value = 1
# row flagged for manual hospital-probability review -> reuse existing format from A136
$wsUnits.Range("A136").Copy()
$wsUnits.Range("A149").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsUnits.Cells.Item(150, 1).Value = "Hb-14D GD Sjúkraþjálfun H"
$wsUnits.Cells.Item(150, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(150, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(150, 4).Value = "home"
$wsUnits.Cells.Item(150, 5).Value = 1

$wsUnits.Cells.Item(151, 1).Value = "La71 GD Laugarásinn meðferðargeðdeild"
$wsUnits.Cells.Item(151, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(151, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(151, 4).Value = "home"
$wsUnits.Cells.Item(151, 5).Value = 1

$wsUnits.Cells.Item(152, 1).Value = "Fv-B1 GD Iðjuþjálfun F"
$wsUnits.Cells.Item(152, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(152, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(152, 4).Value = "home"
$wsUnits.Cells.Item(152, 5).Value = 1
# row flagged for manual hospital-probability review -> reuse existing format from A136
$wsUnits.Range("A136").Copy()
$wsUnits.Range("A152").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsUnits.Cells.Item(153, 1).Value = "Gr-R3 DD Endurhæfingar"
$wsUnits.Cells.Item(153, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(153, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(153, 4).Value = "home"
$wsUnits.Cells.Item(153, 5).Value = 1
# row flagged for manual hospital-probability review -> reuse existing format from A136
$wsUnits.Range("A136").Copy()
$wsUnits.Range("A153").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsUnits.Cells.Item(154, 1).Value = "Gr-R3 GD Læknar endurhæfingar"
$wsUnits.Cells.Item(154, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(154, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(154, 4).Value = "home"
$wsUnits.Cells.Item(154, 5).Value = 1
# row flagged for manual hospital-probability review -> reuse existing format from A136
$wsUnits.Range("A136").Copy()
$wsUnits.Range("A154").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsUnits.Cells.Item(155, 1).Value = "Ei5 GD Augnlækninga"
$wsUnits.Cells.Item(155, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(155, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(155, 4).Value = "home"
$wsUnits.Cells.Item(155, 5).Value = 1

$wsUnits.Cells.Item(156, 1).Value = "Ei5 GD Innkirtla- og efnaskipta"
$wsUnits.Cells.Item(156, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(156, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(156, 4).Value = "home"
$wsUnits.Cells.Item(156, 5).Value = 1

$wsUnits.Cells.Item(157, 1).Value = "Sjúkrahótel Landspítala (Hb-4)"
$wsUnits.Cells.Item(157, 2).Value = "Legudeild"
$wsUnits.Cells.Item(157, 3).Value = "inpatient_ward"
$wsUnits.Cells.Item(157, 4).Value = "inpatient_ward"
$wsUnits.Cells.Item(157, 5).Value = 2

$wsUnits.Cells.Item(158, 1).Value = "Ei5 GD Gigtarlækninga"
$wsUnits.Cells.Item(158, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(158, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(158, 4).Value = "home"
$wsUnits.Cells.Item(158, 5).Value = 1

$wsUnits.Cells.Item(159, 1).Value = "Gr-R1 GD Sjúkraþjálfun G"
$wsUnits.Cells.Item(159, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(159, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(159, 4).Value = "home"
$wsUnits.Cells.Item(159, 5).Value = 1

$wsUnits.Cells.Item(160, 1).Value = "Fv-B7 DD Alm.lyflækninga"
$wsUnits.Cells.Item(160, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(160, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(160, 4).Value = "home"
$wsUnits.Cells.Item(160, 5).Value = 1

$wsUnits.Cells.Item(161, 1).Value = "Öldrunarlækningadeild (Lk-L5)"
$wsUnits.Cells.Item(161, 2).Value = "Legudeild"
$wsUnits.Cells.Item(161, 3).Value = "inpatient_ward_geriatric"
$wsUnits.Cells.Item(161, 4).Value = "inpatient_ward"
$wsUnits.Cells.Item(161, 5).Value = 2

$wsUnits.Cells.Item(162, 1).Value = "Ei5 GD Erfða- og sameindalæknisfræðideild"
$wsUnits.Cells.Item(162, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(162, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(162, 4).Value = "home"
$wsUnits.Cells.Item(162, 5).Value = 1

$wsUnits.Cells.Item(163, 1).Value = "Ei5 DD Gigtarlækningar"
$wsUnits.Cells.Item(163, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(163, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(163, 4).Value = "home"
$wsUnits.Cells.Item(163, 5).Value = 1

$wsUnits.Cells.Item(164, 1).Value = "Ei5 GD Sálfræðiþjónusta"
$wsUnits.Cells.Item(164, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(164, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(164, 4).Value = "home"
$wsUnits.Cells.Item(164, 5).Value = 1

$wsUnits.Cells.Item(165, 1).Value = "Fv-E4 Æðaþræðing-inngripsröntgen"
$wsUnits.Cells.Item(165, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(165, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(165, 4).Value = "home"
$wsUnits.Cells.Item(165, 5).Value = 1

$wsUnits.Cells.Item(166, 1).Value = "Ei5 GD Brjóstamóttaka"
$wsUnits.Cells.Item(166, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(166, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(166, 4).Value = "home"
$wsUnits.Cells.Item(166, 5).Value = 1

$wsUnits.Cells.Item(167, 1).Value = "Ei5 GD Lýtalækninga"
$wsUnits.Cells.Item(167, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(167, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(167, 4).Value = "home"
$wsUnits.Cells.Item(167, 5).Value = 1

$wsUnits.Cells.Item(168, 1).Value = "Næringarstofa"
$wsUnits.Cells.Item(168, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(168, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(168, 4).Value = "home"
$wsUnits.Cells.Item(168, 5).Value = 1

$wsUnits.Cells.Item(169, 1).Value = "Hs-7h GD Verkjamiðstöð"
$wsUnits.Cells.Item(169, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(169, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(169, 4).Value = "home"
$wsUnits.Cells.Item(169, 5).Value = 1

$wsUnits.Cells.Item(170, 1).Value = "Hb-20E Næringarstofa BH"
$wsUnits.Cells.Item(170, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(170, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(170, 4).Value = "home"
$wsUnits.Cells.Item(170, 5).Value = 1

$wsUnits.Cells.Item(171, 1).Value = "Ei5 GD Krabbameins"
$wsUnits.Cells.Item(171, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(171, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(171, 4).Value = "home"
$wsUnits.Cells.Item(171, 5).Value = 1

$wsUnits.Cells.Item(172, 1).Value = "Ei5 GD Kviðarhols- og brjóstaskurðlækninga"
$wsUnits.Cells.Item(172, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(172, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(172, 4).Value = "home"
$wsUnits.Cells.Item(172, 5).Value = 1

$wsUnits.Cells.Item(173, 1).Value = "Fv-Bb DD Lyflækningaþjónusta"
$wsUnits.Cells.Item(173, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(173, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(173, 4).Value = "home"
$wsUnits.Cells.Item(173, 5).Value = 1

$wsUnits.Cells.Item(174, 1).Value = "Kl-Sk SV Samfélagsgeðteymi"
$wsUnits.Cells.Item(174, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(174, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(174, 4).Value = "home"
$wsUnits.Cells.Item(174, 5).Value = 1

$wsUnits.Cells.Item(175, 1).Value = "Fv-B3 GD Talmeinaþjónusta"
$wsUnits.Cells.Item(175, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(175, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(175, 4).Value = "home"
$wsUnits.Cells.Item(175, 5).Value = 1

$wsUnits.Cells.Item(176, 1).Value = "Hb-10D GD Hjartavísindarannsóknir"
$wsUnits.Cells.Item(176, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(176, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(176, 4).Value = "home"
$wsUnits.Cells.Item(176, 5).Value = 1

$wsUnits.Cells.Item(177, 1).Value = "Db12-0h Iðjuþjálfun geðendurhæfing"
$wsUnits.Cells.Item(177, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(177, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(177, 4).Value = "home"
$wsUnits.Cells.Item(177, 5).Value = 1

$wsUnits.Cells.Item(178, 1).Value = "Kl-H10 DD Iðjuþjálfunar Kleppi"
$wsUnits.Cells.Item(178, 2).Value = "Dagdeild"
$wsUnits.Cells.Item(178, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(178, 4).Value = "home"
$wsUnits.Cells.Item(178, 5).Value = 1

$wsUnits.Cells.Item(179, 1).Value = "Fv-G3 GD SBS - inniliggjandi"
$wsUnits.Cells.Item(179, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(179, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(179, 4).Value = "home"
$wsUnits.Cells.Item(179, 5).Value = 1

$wsUnits.Cells.Item(180, 1).Value = "Hb-20E GD Talmeinaþjónusta BH"
$wsUnits.Cells.Item(180, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(180, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(180, 4).Value = "home"
$wsUnits.Cells.Item(180, 5).Value = 1

$wsUnits.Cells.Item(181, 1).Value = "Hb-21A GD Innskriftir kvenna"
$wsUnits.Cells.Item(181, 2).Value = "Göngudeild"
$wsUnits.Cells.Item(181, 3).Value = "outpatient_clinic"
$wsUnits.Cells.Item(181, 4).Value = "home"
$wsUnits.Cells.Item(181, 5).Value = 1

$excel.CutCopyMode = $false

# ---- Restore view state: isolation sheet selection, then unit sheet ----
# ---- (the last-activated sheet becomes the saved "tabSelected" sheet) ----
$wsIsolation.Activate()
$wsIsolation.Range("A37").Select()

$wsUnits.Activate()
try {
  $excel.ActiveWindow.ScrollRow = 159
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$wsUnits.Range("B176").Select()
